$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the municipality name (was just "1")
$ws.Name = "ხელვაჩაური"

# The subtitle "(by population census results)" in A2, and the placeholder
# spacer cells under the title (B1:B3), are no longer used -- clear them
# completely (value + formatting) so they vanish from the sheet.
$ws.Range("A2").Clear()
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()
$ws.Range("B3").Clear()

# Remove the now-empty row 2, shifting the "area" header and the data rows
# up by one.
$ws.Rows.Item(2).Delete()

# Only the 2014 figures are kept going forward -- drop the 1989 and 2002
# columns, which shifts the 2014 column (with its own formatting) into
# column B.
$ws.Range("B:C").Delete()

# The data-row label no longer needs its left edge (it used to line up with
# the now-removed 1989/2002 columns) -- drop that border.
$ws.Range("A5").Borders.Item(7).LineStyle = 0

# Leave the cursor on the now-blank spacer row, as in the saved file.
$ws.Range("A2").Select() | Out-Null
